# Fix bug: exceeded request in google drive
# Update the date in A1 (Hoja1) by one day, and correct the price values
# in D23:D26 on the "Hoja1" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# A1 holds a date serial (45310 -> 45311), keep existing date formatting/style.
$ws.Range("A1").Value = 45311

# Update prices
$ws.Range("D23").Value = 398
$ws.Range("D24").Value = 398
$ws.Range("D25").Value = 400
$ws.Range("D26").Value = 398
